$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.025.12'
$ws.Range('E2').Value = '  +7.11%  '
$ws.Range('D3').Value = '1.883.72'
$ws.Range('E3').Value = '  +5.72%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '249.56'
$ws.Range('E5').Value = '  +2.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4982'
$ws.Range('E7').Value = '  +1.23%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '45.79'
$ws.Range('E8').Value = '  +8.93%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2867'
$ws.Range('E9').Value = '  +7.16%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06575'
$ws.Range('E10').Value = '  +4.94%  '
$ws.Range('D11').Value = '1.888.39'
$ws.Range('E11').Value = '  +6.01%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '17.18'
$ws.Range('E12').Value = '  +4.66%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07212'
$ws.Range('E13').Value = '  +2.56%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6658'
$ws.Range('E14').Value = '  +6.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '85.35'
$ws.Range('E15').Value = '  +6.43%  '
$ws.Range('E16').Value = '  +3.48%  '
$ws.Range('D17').Value = '30.012.42'
$ws.Range('E17').Value = '  +7.24%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.000'
$ws.Range('E18').Value = '  +0.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.91'
$ws.Range('E19').Value = '  +7.85%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007521'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.000'
$ws.Range('E21').Value = '  +0.19%  '
$ws.Range('D22').Value = '2.122.03'
$ws.Range('E22').Value = '  +5.45%  '
$ws.Range('E23').Value = '  +4.51%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.535'
$ws.Range('E24').Value = '  +5.81%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.037'
$ws.Range('E25').Value = '  +3.80%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '144.42'
$ws.Range('E26').Value = '  +1.65%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '134.85'
$ws.Range('E27').Value = '  +23.26%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.974'
$ws.Range('E29').Value = '  +6.03%  '
$ws.Range('E30').Value = '  +0.57%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.170'
$ws.Range('E31').Value = '  -0.70%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.08603'
$ws.Range('E32').Value = '  +3.86%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.883'
$ws.Range('E33').Value = '  +2.31%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05092'
$ws.Range('E34').Value = '  +4.11%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.142'
$ws.Range('E35').Value = '  +6.23%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6879'
$ws.Range('E36').Value = '  +5.32%  '
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.710'
$ws.Range('E37').Value = '  +3.75%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.301'
$ws.Range('E38').Value = '  +12.04%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.757'
$ws.Range('E39').Value = '  +7.01%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9622'
$ws.Range('E40').Value = '  +1.45%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.01632'
$ws.Range('E41').Value = '  +5.01%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.106'
$ws.Range('E42').Value = '  +1.89%  '
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.001'
$ws.Range('E43').Value = '  +0.07%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '103.79'
$ws.Range('E44').Value = '  +3.80%  '
$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4208'
$ws.Range('E45').Value = '  +5.42%  '
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.464'
$ws.Range('E46').Value = '  +3.87%  '
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1255'
$ws.Range('E47').Value = '  +4.40%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05636'
$ws.Range('E48').Value = '  +3.95%  '
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '32.52'
$ws.Range('E49').Value = '  +5.96%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.269'
$ws.Range('E50').Value = '  +2.97%  '
$ws.Range('B51').Value = 'Decentraland'
$ws.Range('C51').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3731'
$ws.Range('E51').Value = '  +7.40%  '
